$d = $word.ActiveDocument

# Locate the paragraph ending the "Things I Learned" list (the
# Input.GetAxis / Input.GetAxisRaw bullet) and append a new list item
# after it containing the new bullet text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Input.GetAxisRaw instead*") {
        $target = $p
    }
}

if ($target -eq $null) {
    $target = $d.Paragraphs.Last
}

$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara = $target.Next()
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.Text = "Very useful Collision Action Matrix of Documentation."
